$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 ({d.tickets[i].*}) — append :formatN() to the numeric ticket-count columns
$ws.Range("C2").Value = "{d.tickets[i].processing:formatN()}"
$ws.Range("D2").Value = "{d.tickets[i].completed:formatN()}"
$ws.Range("E2").Value = "{d.tickets[i].canceled:formatN()}"
$ws.Range("F2").Value = "{d.tickets[i].deferred:formatN()}"
$ws.Range("G2").Value = "{d.tickets[i].closed:formatN()}"
$ws.Range("H2").Value = "{d.tickets[i].new_or_reopened:formatN()}"

# Row 3 ({d.tickets[i+1].*}) — same treatment
$ws.Range("C3").Value = "{d.tickets[i+1].processing:formatN()}"
$ws.Range("D3").Value = "{d.tickets[i+1].completed:formatN()}"
$ws.Range("E3").Value = "{d.tickets[i+1].canceled:formatN()}"
$ws.Range("F3").Value = "{d.tickets[i+1].deferred:formatN()}"
$ws.Range("G3").Value = "{d.tickets[i+1].closed:formatN()}"
$ws.Range("H3").Value = "{d.tickets[i+1].new_or_reopened:formatN()}"

# Apply a numeric ("0") number format to those same cells so the formatted
# value renders as a plain integer.
$ws.Range("C2:H2").NumberFormat = "0"
$ws.Range("C3:H3").NumberFormat = "0"
